# "fixed change depth, added initialize to KPro_2016 for first point on KPRO Start"
#
# The "Sample Sheet" (second worksheet) holds a two-column lookup table:
# column A = friendly field name, column B = the worksheet cell the field
# maps to. Three new auxiliary rows are appended below the existing table
# (which ran through row 19): Aux1/Aux2 (U14/V14) and a third Aux3 (W14)
# row added right after, to fix the change-depth / initialization lookups.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Sheet")
$ws.Activate()

# Aux1 / Aux2 labels first, then their cell references (matches the order
# the fields were actually typed in), and finally the separately-added Aux3.
$ws.Range("A20").Value = "Aux1"
$ws.Range("A21").Value = "Aux2"
$ws.Range("B20").Value = "U14"
$ws.Range("B21").Value = "V14"
$ws.Range("A22").Value = "Aux3"
$ws.Range("B22").Value = "W14"

# Scroll the sheet down a bit and leave the selection on the new Aux2 row.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
